$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = "'47.493.74"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Formula = "'2.492.40"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Formula = "'321.85"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').Formula = "'109.30"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.89%  '
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('D8').Formula = "'1.00"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Formula = "'0.543"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.38%  '
$ws.Range('D10').Formula = "'39.43"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.48%  '
$ws.Range('D11').Formula = "'0.0811"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.31%  '
$ws.Range('E12').Value = '  +0.93%  '
$ws.Range('D13').Formula = "'18.63"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.17%  '
$ws.Range('D14').Formula = "'7.20"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').Formula = "'2.878.81"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.15%  '
$ws.Range('D16').Formula = "'2.490.30"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('D17').Formula = "'0.848"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('D18').Formula = "'47.350.36"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').Formula = "'13.47"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.40%  '
$ws.Range('D20').Formula = "'6.64"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.39%  '
$ws.Range('D21').Formula = "'0.0₃0942"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.89%  '
$ws.Range('D22').Formula = "'2.77"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +16.68%  '
$ws.Range('D23').Formula = "'70.71"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').Formula = "'246.94"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.35%  '
$ws.Range('D25').Formula = "'2.56"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -1.40%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Formula = "'9.99"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Formula = "'2.21"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('D30').Formula = "'0.139"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.42%  '
$ws.Range('D31').Formula = "'34.78"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.46%  '
$ws.Range('D32').Formula = "'49.93"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.95%  '
$ws.Range('D33').Formula = "'20.46"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.22%  '
$ws.Range('D34').Formula = "'5.32"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.31%  '
$ws.Range('D35').Formula = "'0.0788"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.08%  '
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').Formula = "'4.74"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.46%  '
$ws.Range('E38').Value = '  +1.32%  '
$ws.Range('D39').Formula = "'2.93"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.64%  '
$ws.Range('E40').Value = '  +0.49%  '
$ws.Range('D41').Formula = "'22.53"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.66%  '
$ws.Range('E42').Value = '  -1.94%  '
$ws.Range('E43').Value = '  -1.96%  '
$ws.Range('D44').Formula = "'0.0297"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.40%  '
$ws.Range('D45').Formula = "'1.995.58"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.82%  '
$ws.Range('D46').Formula = "'3.04"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.54%  '
$ws.Range('D47').Formula = "'2.05"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.59%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Formula = "'1.78"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.68%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Formula = "'9.09"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.08%  '
$ws.Range('D50').Formula = "'5.22"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.50%  '
$ws.Range('D51').Formula = "'56.82"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.80%  '
